$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the CreatedAt timestamp in A1
$ws.Range("A1").Value = "CreatedAt: 2025-05-28T15:08:33"

# Update Intertie LMP / Energy Loss / Energy Congestion values for hours 17-24 (columns S:Z)
# across the affected data rows, reflecting refreshed IESO report figures.

$rng = $ws.Range("S4:Z4")
$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 70.72
$arr[0,1] = 54.71
$arr[0,2] = 43.99
$arr[0,3] = 44.37
$arr[0,4] = 32.15
$arr[0,5] = 30.76
$arr[0,6] = 14.77
$arr[0,7] = 14.8
$rng.Value = $arr

$rng = $ws.Range("S6:Z6")
$arr = New-Object 'object[,]' 1,8
$arr[0,0] = -0.85
$arr[0,1] = -0.6
$arr[0,2] = -0.57
$arr[0,3] = -0.75
$arr[0,4] = -0.55
$arr[0,5] = -0.22
$arr[0,6] = 0
$arr[0,7] = 0.03
$rng.Value = $arr

$rng = $ws.Range("S9:Z9")
$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 67.39
$arr[0,1] = 52.83
$arr[0,2] = 42.97
$arr[0,3] = 43.14
$arr[0,4] = 31.4
$arr[0,5] = 30.6
$arr[0,6] = 15.11
$arr[0,7] = 15.13
$rng.Value = $arr

$rng = $ws.Range("S11:Z11")
$arr = New-Object 'object[,]' 1,8
$arr[0,0] = -4.18
$arr[0,1] = -2.48
$arr[0,2] = -1.59
$arr[0,3] = -1.98
$arr[0,4] = -1.29
$arr[0,5] = -0.37
$arr[0,6] = 0.33
$arr[0,7] = 0.36
$rng.Value = $arr

$rng = $ws.Range("S14:Z14")
$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 68.27
$arr[0,1] = 52.83
$arr[0,2] = 42.93
$arr[0,3] = 43.14
$arr[0,4] = 31.4
$arr[0,5] = 30.6
$arr[0,6] = 15.11
$arr[0,7] = 15.13
$rng.Value = $arr

$rng = $ws.Range("S15:Z15")
$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 0.9399999999999999
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$rng.Value = $arr

$rng = $ws.Range("S16:Z16")
$arr = New-Object 'object[,]' 1,8
$arr[0,0] = -4.24
$arr[0,1] = -2.48
$arr[0,2] = -1.63
$arr[0,3] = -1.98
$arr[0,4] = -1.29
$arr[0,5] = -0.37
$arr[0,6] = 0.33
$arr[0,7] = 0.36
$rng.Value = $arr

$rng = $ws.Range("S19:Z19")
$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 70.79000000000001
$arr[0,1] = 54.87
$arr[0,2] = 44.16
$arr[0,3] = 44.45
$arr[0,4] = 32.24
$arr[0,5] = 31
$arr[0,6] = 14.94
$arr[0,7] = 14.96
$rng.Value = $arr

$rng = $ws.Range("S21:Z21")
$arr = New-Object 'object[,]' 1,8
$arr[0,0] = -0.78
$arr[0,1] = -0.44
$arr[0,2] = -0.4
$arr[0,3] = -0.67
$arr[0,4] = -0.45
$arr[0,5] = 0.03
$arr[0,6] = 0.16
$arr[0,7] = 0.19
$rng.Value = $arr

$rng = $ws.Range("S24:Z24")
$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 70.79000000000001
$arr[0,1] = 54.87
$arr[0,2] = 44.16
$arr[0,3] = 44.45
$arr[0,4] = 32.24
$arr[0,5] = 31
$arr[0,6] = 14.94
$arr[0,7] = 14.96
$rng.Value = $arr

$rng = $ws.Range("S26:Z26")
$arr = New-Object 'object[,]' 1,8
$arr[0,0] = -0.78
$arr[0,1] = -0.44
$arr[0,2] = -0.4
$arr[0,3] = -0.67
$arr[0,4] = -0.45
$arr[0,5] = 0.03
$arr[0,6] = 0.16
$arr[0,7] = 0.19
$rng.Value = $arr

$rng = $ws.Range("S29:Z29")
$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 71.14
$arr[0,1] = 55.31
$arr[0,2] = 44.47
$arr[0,3] = 44.72
$arr[0,4] = 32.5
$arr[0,5] = 31.38
$arr[0,6] = 15.14
$arr[0,7] = 15.16
$rng.Value = $arr

$rng = $ws.Range("S31:Z31")
$arr = New-Object 'object[,]' 1,8
$arr[0,0] = -0.43
$arr[0,1] = 0
$arr[0,2] = -0.09
$arr[0,3] = -0.4
$arr[0,4] = -0.19
$arr[0,5] = 0.41
$arr[0,6] = 0.36
$arr[0,7] = 0.39
$rng.Value = $arr

$rng = $ws.Range("S34:Z34")
$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 65.76000000000001
$arr[0,1] = 51.07
$arr[0,2] = 41.96
$arr[0,3] = 41.97
$arr[0,4] = 30.58
$arr[0,5] = 30.27
$arr[0,6] = 15.25
$arr[0,7] = 15.26
$rng.Value = $arr

$rng = $ws.Range("S35:Z35")
$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 0.9399999999999999
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$rng.Value = $arr

$rng = $ws.Range("S36:Z36")
$arr = New-Object 'object[,]' 1,8
$arr[0,0] = -6.74
$arr[0,1] = -4.24
$arr[0,2] = -2.6
$arr[0,3] = -3.15
$arr[0,4] = -2.11
$arr[0,5] = -0.7
$arr[0,6] = 0.47
$arr[0,7] = 0.49
$rng.Value = $arr

$rng = $ws.Range("S39:Z39")
$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 70.72
$arr[0,1] = 54.71
$arr[0,2] = 43.99
$arr[0,3] = 44.37
$arr[0,4] = 32.15
$arr[0,5] = 30.76
$arr[0,6] = 14.77
$arr[0,7] = 14.8
$rng.Value = $arr

$rng = $ws.Range("S41:Z41")
$arr = New-Object 'object[,]' 1,8
$arr[0,0] = -0.85
$arr[0,1] = -0.6
$arr[0,2] = -0.57
$arr[0,3] = -0.75
$arr[0,4] = -0.55
$arr[0,5] = -0.22
$arr[0,6] = 0
$arr[0,7] = 0.03
$rng.Value = $arr

$rng = $ws.Range("S44:Z44")
$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 72.36
$arr[0,1] = 55.98
$arr[0,2] = 45.1
$arr[0,3] = 45.62
$arr[0,4] = 32.99
$arr[0,5] = 31.13
$arr[0,6] = 14.79
$arr[0,7] = 14.8
$rng.Value = $arr

$rng = $ws.Range("S46:Z46")
$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 0.8
$arr[0,1] = 0.67
$arr[0,2] = 0.54
$arr[0,3] = 0.5
$arr[0,4] = 0.3
$arr[0,5] = 0.16
$arr[0,6] = 0.01
$arr[0,7] = 0.03
$rng.Value = $arr

$rng = $ws.Range("S49:Z49")
$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 64.19
$arr[0,1] = 47.76
$arr[0,2] = 38.22
$arr[0,3] = 38.76
$arr[0,4] = 28.09
$arr[0,5] = 26.25
$arr[0,6] = 13.1
$arr[0,7] = 12.83
$rng.Value = $arr

$rng = $ws.Range("S51:Z51")
$arr = New-Object 'object[,]' 1,8
$arr[0,0] = -7.38
$arr[0,1] = -7.55
$arr[0,2] = -6.34
$arr[0,3] = -6.36
$arr[0,4] = -4.61
$arr[0,5] = -4.72
$arr[0,6] = -1.68
$arr[0,7] = -1.94
$rng.Value = $arr

$rng = $ws.Range("S54:Z54")
$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 59.89
$arr[0,1] = 46.6
$arr[0,2] = 37.13
$arr[0,3] = 38.7
$arr[0,4] = 28.16
$arr[0,5] = 26.63
$arr[0,6] = 12.69
$arr[0,7] = 12.65
$rng.Value = $arr

$rng = $ws.Range("S56:Z56")
$arr = New-Object 'object[,]' 1,8
$arr[0,0] = -11.68
$arr[0,1] = -8.710000000000001
$arr[0,2] = -7.43
$arr[0,3] = -6.42
$arr[0,4] = -4.53
$arr[0,5] = -4.34
$arr[0,6] = -2.08
$arr[0,7] = -2.11
$rng.Value = $arr

$rng = $ws.Range("S59:Z59")
$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 74.55
$arr[0,1] = 57.49
$arr[0,2] = 46.32
$arr[0,3] = 46.85
$arr[0,4] = 33.91
$arr[0,5] = 31.9
$arr[0,6] = 15.14
$arr[0,7] = 15.12
$rng.Value = $arr

$rng = $ws.Range("S61:Z61")
$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 2.98
$arr[0,1] = 2.18
$arr[0,2] = 1.76
$arr[0,3] = 1.73
$arr[0,4] = 1.22
$arr[0,5] = 0.92
$arr[0,6] = 0.36
$arr[0,7] = 0.35
$rng.Value = $arr

$rng = $ws.Range("S64:Z64")
$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 75.73
$arr[0,1] = 58.41
$arr[0,2] = 47.06
$arr[0,3] = 47.65
$arr[0,4] = 34.45
$arr[0,5] = 32.36
$arr[0,6] = 15.33
$arr[0,7] = 15.3
$rng.Value = $arr

$rng = $ws.Range("S66:Z66")
$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 4.17
$arr[0,1] = 3.1
$arr[0,2] = 2.49
$arr[0,3] = 2.53
$arr[0,4] = 1.76
$arr[0,5] = 1.39
$arr[0,6] = 0.55
$arr[0,7] = 0.54
$rng.Value = $arr

$rng = $ws.Range("S69:Z69")
$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 76.54000000000001
$arr[0,1] = 59.03
$arr[0,2] = 47.56
$arr[0,3] = 48.15
$arr[0,4] = 34.82
$arr[0,5] = 32.67
$arr[0,6] = 15.47
$arr[0,7] = 15.45
$rng.Value = $arr

$rng = $ws.Range("S71:Z71")
$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 4.98
$arr[0,1] = 3.72
$arr[0,2] = 3
$arr[0,3] = 3.03
$arr[0,4] = 2.12
$arr[0,5] = 1.7
$arr[0,6] = 0.7
$arr[0,7] = 0.68
$rng.Value = $arr

$rng = $ws.Range("S74:Z74")
$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 73.09999999999999
$arr[0,1] = 56.32
$arr[0,2] = 45.38
$arr[0,3] = 45.95
$arr[0,4] = 33.22
$arr[0,5] = 31.22
$arr[0,6] = 14.82
$arr[0,7] = 14.8
$rng.Value = $arr

$rng = $ws.Range("S76:Z76")
$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 1.54
$arr[0,1] = 1.01
$arr[0,2] = 0.82
$arr[0,3] = 0.83
$arr[0,4] = 0.53
$arr[0,5] = 0.25
$arr[0,6] = 0.04
$arr[0,7] = 0.03
$rng.Value = $arr

$rng = $ws.Range("S79:Z79")
$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 71.56999999999999
$arr[0,1] = 55.31
$arr[0,2] = 44.56
$arr[0,3] = 45.12
$arr[0,4] = 32.69
$arr[0,5] = 30.97
$arr[0,6] = 14.77
$arr[0,7] = 14.77
$rng.Value = $arr

$rng = $ws.Range("S84:Z84")
$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 59.89
$arr[0,1] = 46.6
$arr[0,2] = 37.13
$arr[0,3] = 38.73
$arr[0,4] = 28.16
$arr[0,5] = 26.65
$arr[0,6] = 12.7
$arr[0,7] = 12.65
$rng.Value = $arr

$rng = $ws.Range("S86:Z86")
$arr = New-Object 'object[,]' 1,8
$arr[0,0] = -11.68
$arr[0,1] = -8.710000000000001
$arr[0,2] = -7.43
$arr[0,3] = -6.39
$arr[0,4] = -4.53
$arr[0,5] = -4.32
$arr[0,6] = -2.07
$arr[0,7] = -2.11
$rng.Value = $arr

$rng = $ws.Range("S89:Z89")
$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 71.14
$arr[0,1] = 55.31
$arr[0,2] = 44.47
$arr[0,3] = 44.72
$arr[0,4] = 32.5
$arr[0,5] = 31.38
$arr[0,6] = 15.14
$arr[0,7] = 15.16
$rng.Value = $arr

$rng = $ws.Range("S91:Z91")
$arr = New-Object 'object[,]' 1,8
$arr[0,0] = -0.43
$arr[0,1] = 0
$arr[0,2] = -0.09
$arr[0,3] = -0.4
$arr[0,4] = -0.19
$arr[0,5] = 0.41
$arr[0,6] = 0.36
$arr[0,7] = 0.39
$rng.Value = $arr
